$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.384.91"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "1.593.59"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "210.13"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "1.557.45"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "64.66"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "26.379.87"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("D20").Value = "211.71"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "2.18"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "145.39"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Value = "1.307.86"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "1.11"
$ws.Range("E39").Value = "  -12.12%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "5.62"
$ws.Range("D43").Value = "62.77"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "1.728.83"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "88.17"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.50"
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("D50").Value = "0.0986"
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("E51").Value = "  -1.42%  "
